# Update "想去人数" (F column) counts that changed between scrapes.
# Same underlying events are listed on both the "展览" sheet and the
# "全部类型" sheet, so each value needs to be updated in both places.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 sheet (rows 5,9,10,12,13,14,15,17,18,19)
$wsExhibit.Range("F5").Value  = 304
$wsExhibit.Range("F9").Value  = 547
$wsExhibit.Range("F10").Value = 4
$wsExhibit.Range("F12").Value = 169
$wsExhibit.Range("F13").Value = 13435
$wsExhibit.Range("F14").Value = 175
$wsExhibit.Range("F15").Value = 19
$wsExhibit.Range("F17").Value = 5536
$wsExhibit.Range("F18").Value = 5576
$wsExhibit.Range("F19").Value = 52

# 全部类型 sheet (rows 21,31,32,34,35,36,37,40,41,42) - mirrors the same entries
$wsAll.Range("F21").Value = 304
$wsAll.Range("F31").Value = 547
$wsAll.Range("F32").Value = 4
$wsAll.Range("F34").Value = 169
$wsAll.Range("F35").Value = 13435
$wsAll.Range("F36").Value = 175
$wsAll.Range("F37").Value = 19
$wsAll.Range("F40").Value = 5536
$wsAll.Range("F41").Value = 5576
$wsAll.Range("F42").Value = 52
